$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row4
$ws.Range("L4").Value = 95
$ws.Range("M4").Value = 15

# Row5
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 2.24
$ws.Range("D5").Value = 45
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 55
$ws.Range("G5").Value = 50
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 76
$ws.Range("M5").Value = 4.9
$ws.Range("N5").Value = 70
$ws.Range("O5").Value = 75
$ws.Range("P5").Value = 90
$ws.Range("Q5").Value = 90

# Row6
$ws.Range("B6").Value = 17.5
$ws.Range("C6").Value = 2.5
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 17.5
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 20
$ws.Range("I6").Value = 20
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 20
$ws.Range("L6").Value = 36.5
$ws.Range("M6").Value = 2.29
$ws.Range("N6").Value = 35
$ws.Range("O6").Value = 35
$ws.Range("P6").Value = 40
$ws.Range("Q6").Value = 40

# Row7
$ws.Range("B7").Value = 31.5
$ws.Range("C7").Value = 2.06
$ws.Range("D7").Value = 29
$ws.Range("E7").Value = 33
$ws.Range("F7").Value = 34
$ws.Range("G7").Value = 24
$ws.Range("I7").Value = 24
$ws.Range("J7").Value = 24
$ws.Range("K7").Value = 24
$ws.Range("L7").Value = 34
$ws.Range("N7").Value = 34
$ws.Range("O7").Value = 34
$ws.Range("P7").Value = 34
$ws.Range("Q7").Value = 29

# Row8
$ws.Range("B8").Value = 134
$ws.Range("D8").Value = 134
$ws.Range("E8").Value = 134
$ws.Range("F8").Value = 134
$ws.Range("G8").Value = 92
$ws.Range("I8").Value = 92
$ws.Range("J8").Value = 92
$ws.Range("K8").Value = 92
$ws.Range("L8").Value = 234
$ws.Range("N8").Value = 234
$ws.Range("O8").Value = 234
$ws.Range("P8").Value = 234
$ws.Range("Q8").Value = 224

# Row9
$ws.Range("B9").Value = 111.2
$ws.Range("C9").Value = 5.84
$ws.Range("D9").Value = 105
$ws.Range("E9").Value = 110
$ws.Range("F9").Value = 118
$ws.Range("G9").Value = 72
$ws.Range("I9").Value = 72
$ws.Range("J9").Value = 72
$ws.Range("K9").Value = 72
$ws.Range("L9").Value = 121.2
$ws.Range("M9").Value = 4.66
$ws.Range("N9").Value = 118
$ws.Range("O9").Value = 118
$ws.Range("P9").Value = 130
$ws.Range("Q9").Value = 130

# Row10
$ws.Range("B10").Value = 106.2
$ws.Range("C10").Value = 8.03
$ws.Range("D10").Value = 96
$ws.Range("E10").Value = 109
$ws.Range("F10").Value = 115
$ws.Range("G10").Value = 57
$ws.Range("I10").Value = 57
$ws.Range("J10").Value = 57
$ws.Range("K10").Value = 57
$ws.Range("L10").Value = 114.5
$ws.Range("M10").Value = 1.8
$ws.Range("N10").Value = 111
$ws.Range("O10").Value = 115
$ws.Range("P10").Value = 117
$ws.Range("Q10").Value = 118

# Row11
$ws.Range("B11").Value = 27
$ws.Range("D11").Value = 27
$ws.Range("E11").Value = 27
$ws.Range("F11").Value = 27
$ws.Range("G11").Value = 24
$ws.Range("I11").Value = 24
$ws.Range("J11").Value = 24
$ws.Range("K11").Value = 24
$ws.Range("L11").Value = 100
$ws.Range("N11").Value = 100
$ws.Range("O11").Value = 100
$ws.Range("P11").Value = 100
$ws.Range("Q11").Value = 100

# Row12
$ws.Range("B12").Value = 80
$ws.Range("D12").Value = 80
$ws.Range("E12").Value = 80
$ws.Range("F12").Value = 80
$ws.Range("G12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 190
$ws.Range("N12").Value = 190
$ws.Range("O12").Value = 190
$ws.Range("P12").Value = 190
$ws.Range("Q12").Value = 190

# Row13
$ws.Range("B13").Value = 134
$ws.Range("C13").Value = 6.63
$ws.Range("D13").Value = 120
$ws.Range("E13").Value = 135
$ws.Range("F13").Value = 140
$ws.Range("G13").Value = 110
$ws.Range("I13").Value = 110
$ws.Range("J13").Value = 110
$ws.Range("K13").Value = 110
$ws.Range("L13").Value = 145
$ws.Range("M13").Value = 6.71
$ws.Range("N13").Value = 140
$ws.Range("O13").Value = 140
$ws.Range("P13").Value = 160
$ws.Range("Q13").Value = 210

